$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value = 5451
$sheet1.Range("F7").Value = 634
$sheet1.Range("F9").Value = 1064
$sheet1.Range("F11").Value = 1510
$sheet1.Range("F12").Value = 4719
$sheet1.Range("F14").Value = 210
$sheet1.Range("F15").Value = 185
$sheet1.Range("F17").Value = 3601
$sheet1.Range("F18").Value = 189
$sheet1.Range("F23").Value = 36
$sheet1.Range("F24").Value = 144
$sheet1.Range("F25").Value = 51
$sheet1.Range("F27").Value = 77
$sheet1.Range("F32").Value = 35

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F4").Value = 5451
$sheet4.Range("F8").Value = 634
$sheet4.Range("F10").Value = 1064
$sheet4.Range("F12").Value = 1510
$sheet4.Range("F13").Value = 4719
$sheet4.Range("F15").Value = 210
$sheet4.Range("F16").Value = 185
$sheet4.Range("F18").Value = 3601
$sheet4.Range("F19").Value = 189
$sheet4.Range("F24").Value = 36
$sheet4.Range("F25").Value = 144
$sheet4.Range("F26").Value = 51
$sheet4.Range("F28").Value = 77
$sheet4.Range("F33").Value = 35
